# Scheduled data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) on the per-job Sheets with freshly pulled market-board data.
# Each block below corresponds to one Leve row (identified by its sheet + row).
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# row 40
$ws.Range("H40").Value = 1795.1666
$ws.Range("I40").Value = 1657.2858
$ws.Range("J40").Value = 1882.909
$ws.Range("K40").Value = 1657.2858
$ws.Range("L40").Value = 1882.909
$ws.Range("M40").Value = -1482.2858
$ws.Range("N40").Value = -2232.909

# row 74
$ws.Range("H74").Value = 2526.2954
$ws.Range("I74").Value = 2154.76
$ws.Range("J74").Value = 3015.158
$ws.Range("K74").Value = 2154.76
$ws.Range("L74").Value = 3015.158
$ws.Range("M74").Value = -1218.76
$ws.Range("N74").Value = -4887.157999999999

# row 76
$ws.Range("H76").Value = 7259.615
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 8034.091
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 8034.091
$ws.Range("M76").Value = -2685
$ws.Range("N76").Value = -8664.091

# row 77
$ws.Range("H77").Value = 2526.2954
$ws.Range("I77").Value = 2154.76
$ws.Range("J77").Value = 3015.158
$ws.Range("K77").Value = 10773.8
$ws.Range("L77").Value = 15075.79
$ws.Range("M77").Value = -6093.800000000001
$ws.Range("N77").Value = -24435.79

# row 79
$ws.Range("H79").Value = 7259.615
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 8034.091
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 8034.091
$ws.Range("M79").Value = -1908
$ws.Range("N79").Value = -10218.091

# row 87
$ws.Range("H87").Value = 39200
$ws.Range("J87").Value = 39200
$ws.Range("L87").Value = 39200
$ws.Range("N87").Value = -41696

# row 90
$ws.Range("H90").Value = 39200
$ws.Range("J90").Value = 39200
$ws.Range("L90").Value = 117600
$ws.Range("N90").Value = -130080

# row 141
$ws.Range("H141").Value = 4999.7554
$ws.Range("I141").Value = 2885.7058
$ws.Range("J141").Value = 11534.091
$ws.Range("K141").Value = 8657.117400000001
$ws.Range("L141").Value = 34602.273
$ws.Range("M141").Value = -3477.117400000001
$ws.Range("N141").Value = -44962.273

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# row 43
$ws.Range("H43").Value = 10710.333
$ws.Range("J43").Value = 10710.333
$ws.Range("L43").Value = 10710.333
$ws.Range("N43").Value = -11336.333

# row 45
$ws.Range("H45").Value = 1912.0667
$ws.Range("I45").Value = 1528.75
$ws.Range("K45").Value = 1528.75
$ws.Range("M45").Value = -1151.75

# row 63
$ws.Range("H63").Value = 1573.1
$ws.Range("I63").Value = 1653.875
$ws.Range("J63").Value = 1250
$ws.Range("K63").Value = 1653.875
$ws.Range("L63").Value = 1250
$ws.Range("M63").Value = -967.875
$ws.Range("N63").Value = -2622

# row 66
$ws.Range("H66").Value = 1573.1
$ws.Range("I66").Value = 1653.875
$ws.Range("J66").Value = 1250
$ws.Range("K66").Value = 8269.375
$ws.Range("L66").Value = 6250
$ws.Range("M66").Value = -4837.375
$ws.Range("N66").Value = -13114

# row 112
$ws.Range("H112").Value = 22189.666
$ws.Range("J112").Value = 22189.666
$ws.Range("L112").Value = 22189.666
$ws.Range("N112").Value = -25143.666

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 21686
$ws.Range("I20").Value = 1089.5
$ws.Range("J20").Value = 44969
$ws.Range("K20").Value = 1089.5
$ws.Range("L20").Value = 44969
$ws.Range("M20").Value = -842.5
$ws.Range("N20").Value = -45463

# row 22
$ws.Range("H22").Value = 175
$ws.Range("I22").Value = 175
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 175
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2
$ws.Range("N22").ClearContents()

# row 80
$ws.Range("H80").Value = 389.38095
$ws.Range("I80").Value = 305.2857
$ws.Range("J80").Value = 431.42856
$ws.Range("K80").Value = 305.2857
$ws.Range("L80").Value = 431.42856
$ws.Range("M80").Value = 692.7143
$ws.Range("N80").Value = -2427.42856

# row 83
$ws.Range("H83").Value = 389.38095
$ws.Range("I83").Value = 305.2857
$ws.Range("J83").Value = 431.42856
$ws.Range("K83").Value = 1526.4285
$ws.Range("L83").Value = 2157.1428
$ws.Range("M83").Value = 3465.5715
$ws.Range("N83").Value = -12141.1428

# row 86
$ws.Range("H86").Value = 1986.9166
$ws.Range("I86").Value = 1153.5454
$ws.Range("K86").Value = 1153.5454
$ws.Range("M86").Value = -30.54539999999997

# row 89
$ws.Range("H89").Value = 1986.9166
$ws.Range("I89").Value = 1153.5454
$ws.Range("K89").Value = 5767.727
$ws.Range("M89").Value = -151.7269999999999

# row 105
$ws.Range("H105").Value = 1876.8182
$ws.Range("I105").Value = 1220
$ws.Range("J105").Value = 2665
$ws.Range("K105").Value = 1220
$ws.Range("L105").Value = 2665
$ws.Range("M105").Value = 527
$ws.Range("N105").Value = -6159

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 5354.846
$ws.Range("I16").Value = 3602.1667
$ws.Range("J16").Value = 6857.143
$ws.Range("K16").Value = 3602.1667
$ws.Range("L16").Value = 6857.143
$ws.Range("M16").Value = -3315.1667
$ws.Range("N16").Value = -7431.143

# row 62
$ws.Range("H62").Value = 4283.7812
$ws.Range("I62").Value = 5099.2856
$ws.Range("J62").Value = 2726.9092
$ws.Range("K62").Value = 5099.2856
$ws.Range("L62").Value = 2726.9092
$ws.Range("M62").Value = -4475.2856
$ws.Range("N62").Value = -3974.9092

# row 65
$ws.Range("H65").Value = 4283.7812
$ws.Range("I65").Value = 5099.2856
$ws.Range("J65").Value = 2726.9092
$ws.Range("K65").Value = 25496.428
$ws.Range("L65").Value = 13634.546
$ws.Range("M65").Value = -22376.428
$ws.Range("N65").Value = -19874.546

# row 113
$ws.Range("H113").Value = 5354.846
$ws.Range("I113").Value = 3602.1667
$ws.Range("J113").Value = 6857.143
$ws.Range("K113").Value = 3602.1667
$ws.Range("L113").Value = 6857.143
$ws.Range("M113").Value = -1432.1667
$ws.Range("N113").Value = -11197.143

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# row 108
$ws.Range("H108").Value = 413.375
$ws.Range("I108").Value = 413.375
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 1240.125
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 1639.875
$ws.Range("N108").ClearContents()

# row 109
$ws.Range("H109").Value = 556.5
$ws.Range("I109").Value = 556.5
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1669.5
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -629.5
$ws.Range("N109").ClearContents()

# row 110
$ws.Range("H110").Value = 2013
$ws.Range("I110").Value = 500
$ws.Range("K110").Value = 1500
$ws.Range("M110").Value = 2590

# row 113
$ws.Range("H113").Value = 656826.5600000001
$ws.Range("I113").Value = 448.9375
$ws.Range("J113").Value = 1181928.8
$ws.Range("K113").Value = 1346.8125
$ws.Range("L113").Value = 3545786.4
$ws.Range("M113").Value = 823.1875
$ws.Range("N113").Value = -3550126.4

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# row 80
$ws.Range("H80").Value = 2846.4614
$ws.Range("I80").Value = 2090
$ws.Range("J80").Value = 4056.8
$ws.Range("K80").Value = 2090
$ws.Range("L80").Value = 4056.8
$ws.Range("M80").Value = -1092
$ws.Range("N80").Value = -6052.8

# row 83
$ws.Range("H83").Value = 2846.4614
$ws.Range("I83").Value = 2090
$ws.Range("J83").Value = 4056.8
$ws.Range("K83").Value = 10450
$ws.Range("L83").Value = 20284
$ws.Range("M83").Value = -5458
$ws.Range("N83").Value = -30268

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 2581.5715
$ws.Range("I46").Value = 2511.8333
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 2511.8333
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -2323.8333
$ws.Range("N46").Value = -3376

# row 64
$ws.Range("H64").Value = 39997.4
$ws.Range("J64").Value = 39997.4
$ws.Range("L64").Value = 39997.4
$ws.Range("N64").Value = -40447.4

# row 67
$ws.Range("H67").Value = 39997.4
$ws.Range("J67").Value = 39997.4
$ws.Range("L67").Value = 39997.4
$ws.Range("N67").Value = -41557.4

# row 110
$ws.Range("H110").Value = 29234.4
$ws.Range("J110").Value = 29234.4
$ws.Range("L110").Value = 29234.4
$ws.Range("N110").Value = -37414.4
